$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (current "Tipo" column) to make room for "MAE"
$ws.Columns("D").Insert()

# New header in D1 (match formatting of the other header cells)
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").Borders.Weight = 2

# Update B (MSE) and C (R2) values, and set D (MAE) values
$ws.Range("B2").Value = 0.4994742218335629
$ws.Range("C2").Value = 0.9900544075492175
$ws.Range("D2").Value = 0.5787430976917186

$ws.Range("B3").Value = 0.2462143248592855
$ws.Range("C3").Value = 0.9951877348287995
$ws.Range("D3").Value = 0.389487354966141

$ws.Range("B4").Value = 0.2813913072590536
$ws.Range("C4").Value = 0.9945866569617249
$ws.Range("D4").Value = 0.4296535512881828

$ws.Range("B5").Value = 0.414215823181585
$ws.Range("C5").Value = 0.9918322310788845
$ws.Range("D5").Value = 0.4964868541284111

$ws.Range("B6").Value = 0.4962020155763345
$ws.Range("C6").Value = 0.985425246969502
$ws.Range("D6").Value = 0.514623945945233

$ws.Range("B7").Value = 0.09566416833399247
$ws.Range("C7").Value = 0.9986778952911759
$ws.Range("D7").Value = 0.2515039606869272

$ws.Range("B8").Value = 0.03715217076155934
$ws.Range("C8").Value = 0.9996170804977886
$ws.Range("D8").Value = 0.1350414551024903

$ws.Range("B9").Value = 0.106471268435573
$ws.Range("C9").Value = 0.9993657810326021
$ws.Range("D9").Value = 0.2232323037118775

$ws.Range("B10").Value = 0.0680855673228315
$ws.Range("C10").Value = 0.9987562379256355
$ws.Range("D10").Value = 0.2061592980343896

$ws.Range("B11").Value = 0.1292650287010617
$ws.Range("C11").Value = 0.990449428770573
$ws.Range("D11").Value = 0.281317959080314

$ws.Range("B12").Value = 0.0516315487418447
$ws.Range("C12").Value = 0.9984816042374429
$ws.Range("D12").Value = 0.1667844346675577

$ws.Range("B13").Value = 0.06039727671964274
$ws.Range("C13").Value = 0.9994271047633952
$ws.Range("D13").Value = 0.1786501536999743

$ws.Range("B14").Value = 0.05580559739108672
$ws.Range("C14").Value = 0.9992444266590271
$ws.Range("D14").Value = 0.189897783702991
